$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -9
    4  = -4
    5  = -5
    7  = -4
    9  = -2
    10 = 4
    14 = -2
    15 = -1
    16 = -13
    18 = -5
    21 = -4
    22 = -4
    24 = -2
    27 = -1
    29 = -3
    30 = -4
    32 = 0
    34 = 6
    35 = -5
    40 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
